# Auto-generated edit script: applies per-cell value updates to the
# Ravana_Profits workbook's per-job sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Each block updates one leve row's market-price-derived columns (H-N).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 222.5
$ws.Range("I12").Value = 195
$ws.Range("K12").Value = 195
$ws.Range("M12").Value = -25
$ws.Range("H15").Value = 1986.7778
$ws.Range("I15").Value = 1986.7778
$ws.Range("K15").Value = 5960.3334
$ws.Range("M15").Value = -5791.3334
$ws.Range("H62").Value = 2332.3333
$ws.Range("I62").Value = 2399
$ws.Range("K62").Value = 2399
$ws.Range("M62").Value = -1775
$ws.Range("H65").Value = 2332.3333
$ws.Range("I65").Value = 2399
$ws.Range("K65").Value = 11995
$ws.Range("M65").Value = -8875
$ws.Range("H80").Value = 594.2
$ws.Range("I80").Value = 602.8570999999999
$ws.Range("J80").Value = 574
$ws.Range("K80").Value = 1808.5713
$ws.Range("L80").Value = 1722
$ws.Range("M80").Value = -810.5712999999998
$ws.Range("N80").Value = -3718
$ws.Range("H83").Value = 594.2
$ws.Range("I83").Value = 602.8570999999999
$ws.Range("J83").Value = 574
$ws.Range("K83").Value = 5425.7139
$ws.Range("L83").Value = 5166
$ws.Range("M83").Value = -433.7138999999997
$ws.Range("N83").Value = -15150
$ws.Range("H86").Value = 4313.3335
$ws.Range("J86").Value = 4222.5
$ws.Range("L86").Value = 4222.5
$ws.Range("N86").Value = -6468.5
$ws.Range("H89").Value = 4313.3335
$ws.Range("J89").Value = 4222.5
$ws.Range("L89").Value = 21112.5
$ws.Range("N89").Value = -32344.5
$ws.Range("H106").Value = 7994.5
$ws.Range("I106").Value = 7989
$ws.Range("J106").Value = 8000
$ws.Range("K106").Value = 7989
$ws.Range("L106").Value = 8000
$ws.Range("M106").Value = -7358
$ws.Range("N106").Value = -9262
$ws.Range("H132").Value = 1031.6
$ws.Range("I132").Value = 1031.6
$ws.Range("K132").Value = 3094.8
$ws.Range("M132").Value = -564.7999999999997
$ws.Range("H138").Value = 3734.2896
$ws.Range("J138").Value = 4043.7
$ws.Range("L138").Value = 12131.1
$ws.Range("N138").Value = -22411.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2696.75
$ws.Range("J4").Value = 2931
$ws.Range("L4").Value = 2931
$ws.Range("N4").Value = -3163
$ws.Range("H32").Value = 5569.5635
$ws.Range("I32").Value = 5656.037
$ws.Range("K32").Value = 5656.037
$ws.Range("M32").Value = -5369.037
$ws.Range("H102").Value = 7449.5
$ws.Range("I102").Value = 7449.5
$ws.Range("K102").Value = 7449.5
$ws.Range("M102").Value = -5827.5
$ws.Range("H122").Value = 8228
$ws.Range("I122").Value = 9637.333000000001
$ws.Range("K122").Value = 28911.999
$ws.Range("M122").Value = -26461.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1219.25
$ws.Range("I94").Value = 917.2857
$ws.Range("K94").Value = 917.2857
$ws.Range("M94").Value = -466.2857
$ws.Range("H107").Value = 739.5
$ws.Range("I107").Value = 739.5
$ws.Range("K107").Value = 739.5
$ws.Range("M107").Value = 1180.5
$ws.Range("H137").Value = 124900
$ws.Range("J137").Value = 124900
$ws.Range("L137").Value = 124900
$ws.Range("N137").Value = -135100
$ws.Range("H138").Value = 124900
$ws.Range("J138").Value = 124900
$ws.Range("L138").Value = 124900
$ws.Range("N138").Value = -135180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1800
$ws.Range("I16").Value = 1800
$ws.Range("K16").Value = 1800
$ws.Range("M16").Value = -1513
$ws.Range("H58").Value = 2954.5386
$ws.Range("I58").Value = 2898
$ws.Range("K58").Value = 2898
$ws.Range("M58").Value = -2695
$ws.Range("H99").Value = 9337.666999999999
$ws.Range("I99").Value = 9999.5
$ws.Range("J99").Value = 8014
$ws.Range("K99").Value = 9999.5
$ws.Range("L99").Value = 8014
$ws.Range("M99").Value = -8501.5
$ws.Range("N99").Value = -11010
$ws.Range("H107").Value = 2174.2354
$ws.Range("I107").Value = 1086.1666
$ws.Range("J107").Value = 2767.7273
$ws.Range("K107").Value = 1086.1666
$ws.Range("L107").Value = 2767.7273
$ws.Range("M107").Value = 833.8334
$ws.Range("N107").Value = -6607.7273
$ws.Range("H113").Value = 1800
$ws.Range("I113").Value = 1800
$ws.Range("K113").Value = 1800
$ws.Range("M113").Value = 370
$ws.Range("H126").Value = 9337.666999999999
$ws.Range("I126").Value = 9999.5
$ws.Range("J126").Value = 8014
$ws.Range("K126").Value = 29998.5
$ws.Range("L126").Value = 24042
$ws.Range("M126").Value = -27528.5
$ws.Range("N126").Value = -28982
$ws.Range("H136").Value = 2954.5386
$ws.Range("I136").Value = 2898
$ws.Range("K136").Value = 8694
$ws.Range("M136").Value = -6144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100.5
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H5").Value = 377.5
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 377.5
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1132.5
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -1356.5
$ws.Range("H18").Value = 3227.1
$ws.Range("I18").Value = 569.75
$ws.Range("K18").Value = 1709.25
$ws.Range("M18").Value = -1540.25
$ws.Range("H34").Value = 4510.3335
$ws.Range("J34").Value = 6082.5
$ws.Range("L34").Value = 18247.5
$ws.Range("N34").Value = -18415.5
$ws.Range("H122").Value = 101647.8
$ws.Range("I122").Value = 1496.2
$ws.Range("K122").Value = 13465.8
$ws.Range("M122").Value = -11015.8
$ws.Range("H135").Value = 377.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 377.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 3397.5
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -8467.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 406.7857
$ws.Range("I97").Value = 446.68182
$ws.Range("K97").Value = 446.68182
$ws.Range("M97").Value = 49.31817999999998
$ws.Range("H122").Value = 2086.6667
$ws.Range("I122").Value = 1130.25
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 3390.75
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -940.75
$ws.Range("N122").Value = -16898.5
$ws.Range("H132").Value = 4437.125
$ws.Range("J132").Value = 4428.143
$ws.Range("L132").Value = 13284.429
$ws.Range("N132").Value = -18344.429
$ws.Range("H136").Value = 84989.5
$ws.Range("J136").Value = 84989.5
$ws.Range("L136").Value = 254968.5
$ws.Range("N136").Value = -260068.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1402
$ws.Range("I7").Value = 1402
$ws.Range("K7").Value = 1402
$ws.Range("M7").Value = -1290
$ws.Range("H40").Value = 2002
$ws.Range("I40").Value = 2002
$ws.Range("K40").Value = 2002
$ws.Range("M40").Value = -1866
$ws.Range("H122").Value = 9599.833000000001
$ws.Range("I122").Value = 6799.5
$ws.Range("K122").Value = 20398.5
$ws.Range("M122").Value = -17948.5
$ws.Range("H126").Value = 1402
$ws.Range("I126").Value = 1402
$ws.Range("K126").Value = 4206
$ws.Range("M126").Value = -1736

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1999
$ws.Range("I96").Value = 1999
$ws.Range("K96").Value = 1999
$ws.Range("M96").Value = -626
$ws.Range("H126").Value = 1631.4166
$ws.Range("I126").Value = 1631.4166
$ws.Range("K126").Value = 4894.2498
$ws.Range("M126").Value = -2424.2498
$ws.Range("H132").Value = 4698.6665
$ws.Range("I132").Value = 4548.5
$ws.Range("K132").Value = 13645.5
$ws.Range("M132").Value = -11115.5
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

